$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("D6").Value = "[Python Pandas] pandas table sorting by other table's column"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Python-Pandas-pandas-table-sorting-by-other-tables-column"

# Row 9
$ws.Range("D9").Value = "데이터 애널리틱스 (Data Analytics) 석사과정 강의목록"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/msda-course-works/#utm_source=rss&utm_medium=rss&utm_campaign=msda-course-works"

# Row 23
$ws.Range("D23").Value = "[ TF Everywhere 행사 영상 및 메이킹 영상 공유]`n안녕하세요! 어제 날짜로 TF Everywhere 텐플마을에 오신것을 환영합니다"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2758"

# Row 37
$ws.Range("D37").Value = "[Paper Review] QANet: Combining Local Convolution  with Global Self-Attention for Reading Comprehension"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1450&mod=document&pageid=1"

# Row 39
$ws.Range("D39").Value = "How to Not Misunderstand Correlation"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/How-to-Not-Misunderstand-Correlation-1"

# Row 51
$ws.Range("D51").Value = "[sqlite3] JOIN으로 서로 다른 테이블의 컬럼들 붙이기(내부 조인, 외부 조인)"
$ws.Range("E51").Value = "https://bskyvision.com/1136"
